# Update absenteeism data rows 2-11 with new values as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 66158
$ws.Range("B2").Value = "Luiz Felipe Correia"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 6675.82

# Row 3
$ws.Range("A3").Value = 88523
$ws.Range("B3").Value = "Luiz Otávio Silveira"
$ws.Range("C3").Value = "Jurídico"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45087
$ws.Range("G3").Value = 5676.45

# Row 4
$ws.Range("A4").Value = 32684
$ws.Range("B4").Value = "Davi Lucca da Paz"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45100
$ws.Range("G4").Value = 5063.57

# Row 5
$ws.Range("A5").Value = 94968
$ws.Range("B5").Value = "Sr. Vitor Hugo Oliveira"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 10223.88

# Row 6
$ws.Range("A6").Value = 90486
$ws.Range("B6").Value = "Alana Jesus"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45096
$ws.Range("G6").Value = 8929.059999999999

# Row 7
$ws.Range("A7").Value = 38736
$ws.Range("B7").Value = "Srta. Julia Sales"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45088
$ws.Range("G7").Value = 10905.71

# Row 8
$ws.Range("A8").Value = 70633
$ws.Range("B8").Value = "Laura Alves"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45083
$ws.Range("G8").Value = 2502.34

# Row 9
$ws.Range("A9").Value = 58941
$ws.Range("B9").Value = "Vitória Ferreira"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 8
$ws.Range("G9").Value = 4992.12

# Row 10
$ws.Range("A10").Value = 77812
$ws.Range("B10").Value = "Srta. Joana Rezende"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 5178.49

# Row 11
$ws.Range("A11").Value = 85837
$ws.Range("B11").Value = "Amanda das Neves"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 5409.13
